$wb = $excel.ActiveWorkbook

# Work on the "train" worksheet
$ws = $wb.Worksheets.Item("train")

# List of hole_id values for rows 2..49 (in order)
$holeIds = @(
    "LBU_02_4",
    "LBU_07_01",
    "LBU_05_05",
    "MHZ_08_01",
    "MHZ_08_03",
    "LBU_05_26",
    "MHZ_12_04",
    "LBU_96_2",
    "LBU_05_12",
    "LBU_05_11",
    "LBU_98_6",
    "LBU_05_07",
    "LBU_05_23",
    "LBU_05_03",
    "MHZ_12_03",
    "LBU_98_1",
    "LBU_05_30",
    "LBU_98_2",
    "MHZ_08_04",
    "LBU_05_01",
    "LBU_96_3",
    "LBU_01_2",
    "LBU_07_02",
    "LBU_01_3",
    "LBU_05_18",
    "LBU_87_3",
    "LBU_05_13",
    "LBU_05_22",
    "LBU_05_10",
    "LBU_87_2",
    "LBU_05_25",
    "LBU_05_16",
    "LBU_05_28",
    "LBU_87_1",
    "MHZ_08_02",
    "LBU_98_7",
    "LBU_87_4",
    "LBU_87_6",
    "LBU_05_06",
    "LBU_05_15",
    "LBU_05_17",
    "LBU_01_1",
    "LBU_05_04",
    "MHZ_08_05",
    "LBU_02_3",
    "MHZ_12_01",
    "LBU_05_20",
    "LBU_87_5"
)

# Add the new "hole_id" header in A1, matching the bold/centered header style
# already used by the other header cells (e.g. B1)
$ws.Range("A1").Value = "hole_id"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Replace numeric index values in A2:A49 with the text hole_id values
for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}

$wb.Save()
